# Update the "取得日時" (retrieved datetime) column for rows 2-7 on the
# ランサーズ sheet to reflect the new append timestamp.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-15 01:48:52"

for ($row = 2; $row -le 7; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
